# Hortaliza, Vega Monumental Concepción - Lechuga.xlsx
# Insert two new weekly price records right after the existing row 634
# (pushing every subsequent record down by two rows) and populate them
# with the new report data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 635-636; everything currently at row 635
# onward shifts down to 637 onward, carrying its formatting with it.
$ws.Range("A635:A636").EntireRow.Insert()

# New row 635
$ws.Range("A635").Value = 11
$ws.Range("B635").Value = "Vega Monumental Concepción"
$ws.Range("C635").Value = "Bíobío"
$ws.Range("D635").Value = 44776
$ws.Range("E635").Value = 8
$ws.Range("F635").Value = 100112033
$ws.Range("G635").Value = "Lechuga"
$ws.Range("H635").Value = "Conconina(o)"
$ws.Range("I635").Value = "Primera"
$ws.Range("J635").Value = 200
$ws.Range("K635").Value = 6500
$ws.Range("L635").Value = 7000
$ws.Range("M635").Value = 6750
$ws.Range("N635").Value = "`$/caja 10 unidades"
$ws.Range("O635").Value = "Región de Valparaíso"
$ws.Range("P635").Value = 675
$ws.Range("Q635").Value = 10
$ws.Range("R635").Value = "Hortaliza"

# New row 636
$ws.Range("A636").Value = 11
$ws.Range("B636").Value = "Vega Monumental Concepción"
$ws.Range("C636").Value = "Bíobío"
$ws.Range("D636").Value = 44776
$ws.Range("E636").Value = 8
$ws.Range("F636").Value = 100112033
$ws.Range("G636").Value = "Lechuga"
$ws.Range("H636").Value = "Escarola"
$ws.Range("I636").Value = "Primera"
$ws.Range("J636").Value = 250
$ws.Range("K636").Value = 9000
$ws.Range("L636").Value = 10000
$ws.Range("M636").Value = 9400
$ws.Range("N636").Value = "`$/caja 15 unidades"
$ws.Range("O636").Value = "Región de Coquimbo"
$ws.Range("P636").Value = 627
$ws.Range("Q636").Value = 15
$ws.Range("R636").Value = "Hortaliza"
